# job_history.xlsx edit:
# - New job run (row 17) filled in with fewer variables (no covariate), mean moved up.
# - 4 new blank separator rows inserted (rows 18-21) with a lighter (borderless) style.
# - What used to be the trailing blank row (old row 18) is now row 22.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Insert 4 new blank rows before the old trailing blank row (old row 18) ---
$ws.Range("A18:A21").EntireRow.Insert()

# --- 2. Fill in row 17 with the new job's data ---
$ws.Range("A17").Value = "ukb51139_subset.csv"
$ws.Range("B17").Value = "28012 x 145"
$ws.Range("C17").Value = "all"
$ws.Range("D17").Value = "no events"
$ws.Range("E17").Value = "> 140/80"
$ws.Range("F17").Value = "zscore"
$ws.Range("G17").Value = "median"
$ws.Range("H17").Value = "none"
$ws.Range("I17").Value = 50
$ws.Range("K17").Value = 37
$ws.Range("L17").Value = "100.3 & 101.3"
$ws.Range("M17").Value = "85.0 & 84.3"
$ws.Range("N17").Value = 17
$ws.Range("O17").Value = 3.51

# --- 3. Restyle the 4 new blank rows (18-21): drop the border, keep the numeric
#        formats / right alignment that the I/K/N/O columns use elsewhere.
#        (Columns A-H,J,L,M already inherited the correct borderless "general"
#        style from the insert above, so only I/K/N/O need adjusting.) ---
$ws.Range("I18:I21").Style = "Normal"
$ws.Range("I18:I21").NumberFormat = "#,##0"
$ws.Range("I18:I21").HorizontalAlignment = -4152

$ws.Range("K18:K21").Style = "Normal"
$ws.Range("K18:K21").NumberFormat = "#,##0"
$ws.Range("K18:K21").HorizontalAlignment = -4152

$ws.Range("N18:N21").Style = "Normal"
$ws.Range("N18:N21").NumberFormat = "#,##0"
$ws.Range("N18:N21").HorizontalAlignment = -4152

$ws.Range("O18:O21").Style = "Normal"
$ws.Range("O18:O21").NumberFormat = "#,##0.00"
$ws.Range("O18:O21").HorizontalAlignment = -4152

# --- 4. Row heights: new separator rows are slightly shorter (18.75) ---
$ws.Rows.Item(18).RowHeight = 18.75
$ws.Rows.Item(19).RowHeight = 18.75
$ws.Rows.Item(20).RowHeight = 18.75
$ws.Rows.Item(21).RowHeight = 18.75
